$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.836.28"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.34"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3885"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.59"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.351"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08473"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.061"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.096"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001316"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.652.14"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.975"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.842.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.435"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.970"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.415"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.857"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.506"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.834.95"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.022"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08204"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.681"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02911"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2678"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09167"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7576"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.425"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6957"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.460"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.104"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08293"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.227"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.50%  "
